# Update the computed answers in the two-digit / one-digit division table.
# Cells are addressed by (row, column) in the single table that holds the
# exercise grid, since several old/new values repeat elsewhere in the
# document and a blind global Find/Replace would be ambiguous/unsafe.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "88÷9=9, 7" },
    @{ Row = 1;  Col = 2; New = "81÷8=10, 1" },
    @{ Row = 1;  Col = 3; New = "16÷8=2, 0" },
    @{ Row = 1;  Col = 4; New = "11÷2=5, 1" },
    @{ Row = 1;  Col = 5; New = "22÷5=4, 2" },

    @{ Row = 5;  Col = 1; New = "51÷4=12, 3" },
    @{ Row = 5;  Col = 2; New = "18÷3=6, 0" },
    @{ Row = 5;  Col = 3; New = "95÷2=47, 1" },
    @{ Row = 5;  Col = 4; New = "32÷6=5, 2" },
    @{ Row = 5;  Col = 5; New = "53÷2=26, 1" },

    @{ Row = 9;  Col = 1; New = "65÷9=7, 2" },
    @{ Row = 9;  Col = 2; New = "15÷9=1, 6" },
    @{ Row = 9;  Col = 3; New = "95÷5=19, 0" },
    @{ Row = 9;  Col = 4; New = "85÷9=9, 4" },
    @{ Row = 9;  Col = 5; New = "30÷6=5, 0" },

    @{ Row = 13; Col = 1; New = "55÷4=13, 3" },
    @{ Row = 13; Col = 2; New = "82÷3=27, 1" },
    @{ Row = 13; Col = 3; New = "62÷5=12, 2" },
    @{ Row = 13; Col = 4; New = "74÷3=24, 2" },
    @{ Row = 13; Col = 5; New = "20÷3=6, 2" },

    @{ Row = 17; Col = 1; New = "65÷9=7, 2" },
    @{ Row = 17; Col = 2; New = "30÷3=10, 0" },
    @{ Row = 17; Col = 3; New = "31÷3=10, 1" },
    @{ Row = 17; Col = 4; New = "49÷7=7, 0" },
    @{ Row = 17; Col = 5; New = "15÷7=2, 1" }
)

foreach ($u in $updates) {
    $cell = $table.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $u.New
}
